$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 6; $row++) {
    $ws.Range("B$row").Value = "sports_club_coed"
    $ws.Range("C$row").Value = "Fun Fit"
}
